# Updating task data a bit.
# This script applies the task-text / status updates to the Gantt chart
# project plan workbook, matching the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Column A (Task) text updates - append completion / collaborator info
# (ordered to match original row order so new shared-string entries are
# appended in the same sequence as the canonical edit)
# ---------------------------------------------------------------------
$ws.Range("A15").Value = ">Understand MovieExplorer Code/Algo (All: Done)"

$ws.Range("A10").Value = ">LaTeX Compile & Submission (Rocko: Done)"
$ws.Range("A11").Value = ">Paper (All, Done)"
$ws.Range("A12").Value = ">Slides (Rocko: Done)"
$ws.Range("A13").Value = ">Video (Daniel: Done)"

$ws.Range("A18").Value = ">Implement Matrix Factor Algo (Jonathan: Done)"
$ws.Range("A19").Value = ">Implement Word2Vec Algo (Daniel, Yi: Done)"

$ws.Range("A21").Value = ">Embed Movies in Taste Space (Daniel: Done)"
$ws.Range("A22").Value = ">Implement TSNE (Rocko, Jonathan, Daniel: Done)"
$ws.Range("A23").Value = ">Pre-cluster movies (Rocko, Jonathan, Daniel: Done)"

# ---------------------------------------------------------------------
# Column B (Status) - mark newly in-progress tasks as WIP
# ---------------------------------------------------------------------
$ws.Range("B26").Value = "WIP"
$ws.Range("B27").Value = "WIP"
$ws.Range("B28").Value = "WIP"

# ---------------------------------------------------------------------
# Column A now has longer text -> widen (best-fit) column A
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 43.8

# ---------------------------------------------------------------------
# Update the selected / active cell on the sheet
# ---------------------------------------------------------------------
$ws.Range("A24").Select() | Out-Null
